# Perubahan nama variabel di jps optimize
# Re-sort the "Waktu Pencarian" data rows (4:65) by Jumlah Kombinasi (A) then
# Kombinasi (B) ascending, and switch the active sheet/selection from
# "Panjang Jalur" (cell G1:G1048576 selected) to "Waktu Pencarian" (cell C7).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Waktu Pencarian")

$sortRange = $ws1.Range("A4:G65")
$key1 = $ws1.Range("A4:A65")
$key2 = $ws1.Range("B4:B65")
$sortRange.Sort($key1, 1, $key2, $null, 1)

$ws1.Activate()
$ws1.Range("C7").Select()
